$d = $word.ActiveDocument

$replacements = @(
    @("2025-02-05 Wednesday", "2025-02-06 Thursday"),
    @("124÷5=", "959÷7="),
    @("759÷8=", "764÷5="),
    @("787÷8=", "408÷4="),
    @("887÷7=", "363÷2="),
    @("101÷7=", "816÷2="),
    @("873÷6=", "573÷6="),
    @("886÷3=", "957÷9="),
    @("769÷2=", "697÷2="),
    @("924÷6=", "578÷2="),
    @("299÷7=", "984÷8="),
    @("407÷6=", "782÷8="),
    @("930÷9=", "893÷5="),
    @("293÷2=", "260÷7="),
    @("134÷6=", "330÷6="),
    @("445÷2=", "565÷4="),
    @("724÷9=", "796÷9="),
    @("772÷9=", "857÷7="),
    @("728÷8=", "265÷3="),
    @("910÷9=", "759÷3="),
    @("988÷7=", "694÷9="),
    @("708÷4=", "377÷8="),
    @("431÷4=", "797÷2="),
    @("432÷9=", "464÷7="),
    @("472÷2=", "347÷2="),
    @("403÷6=", "940÷7=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
